$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 11:52"

# --- Swap province/city labels (column A) between rows 13/14 and 36/37 ---
$ws.Range("A13").Value = "Alacant/Alicante"
$ws.Range("A14").Value = "Araba/Alava"

$ws.Range("A36").Value = "Castello/Castellon"
$ws.Range("A37").Value = "Guadalajara"

# --- Update numeric data cells ---
# Row 9
$ws.Range("B9").Value = 4367
$ws.Range("C9").Value = 1163
$ws.Range("D9").Value = 2825
$ws.Range("E9").Value = 379

# Row 13
$ws.Range("B13").Value = 2962
$ws.Range("C13").Value = 705
$ws.Range("D13").Value = 1939
$ws.Range("E13").Value = 318

# Row 14
$ws.Range("B14").Value = 2868
$ws.Range("C14").Value = 4514
$ws.Range("D14").Value = 4603
$ws.Range("E14").Value = 237

# Row 36
$ws.Range("B36").Value = 997
$ws.Range("C36").Value = 217
$ws.Range("D36").Value = 681
$ws.Range("E36").Value = 99

# Row 37
$ws.Range("B37").Value = 994
$ws.Range("C37").Value = 1766
$ws.Range("D37").Value = 9401
$ws.Range("E37").Value = 134
